$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.213.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.23%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.863.29'
$ws.Range('D3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.18%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.00%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4665'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.50%  '

$ws.Range('E8').Value = '  -0.39%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06535'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.00%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.92%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07855'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.96%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.44'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.28%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.867.52'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.07%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.105'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.03%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6727'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.26%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '280.37'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.23%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.207.31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.28%  '

$ws.Range('E18').Value = '  +0.04%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.519'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.15%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.65'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.17%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.112.32'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.74%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007280'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.29%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9999'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.14%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.148'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.53%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.198'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.67'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.63%  '

$ws.Range('E27').Value = '  -0.49%  '

$ws.Range('E28').Value = '  -3.29%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.378'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.05%  '

$ws.Range('E30').Value = '  -0.62%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.418'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.04%  '

$ws.Range('E32').Value = '  -0.76%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.085'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.42%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04693'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.49%  '

$ws.Range('E35').Value = '  +1.46%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7058'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.19%  '

$ws.Range('E37').Value = '  +0.61%  '

$ws.Range('E38').Value = '  -0.83%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.529'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.12%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.225'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.26%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.26'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.34%  '

$ws.Range('E42').Value = '  -1.63%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8478'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.41%  '

$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4165'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.46%  '

$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '103.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9998'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.04%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.201'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.14%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.151'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.98%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '934.11'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.11%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.09'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.50%  '

$ws.Range('E51').Value = '  -1.98%  '
